$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. First paragraph: pad the existing sentence with two trailing spaces and
#    append a new, differently-formatted (dark red) run containing the
#    "(This is a change ... )" annotation.
# ---------------------------------------------------------------------------

$originalSentence = "This is a Microsoft word document."
$paddedSentence    = "This is a Microsoft word document.  "

$found = $d.Content.Find.Execute(
    $originalSentence, $false, $false, $false, $false, $false,
    $true, 1, $false, $paddedSentence, 2)

$p1 = $d.Paragraphs(1)
$paraRange = $p1.Range

# Collapse to just before the paragraph mark so the new run is appended
# inside the same paragraph (not after it).
$insertPoint = $d.Range($paraRange.End - 1, $paraRange.End - 1)
$insertStart = $insertPoint.End

$dash = [string][char]0x2013
$annotation = "(This is a change " + $dash + " Version for branch alternate)"

$insertPoint.InsertAfter($annotation)

# Re-grab the just-inserted text as its own Range and color it.
$annotationRange = $d.Range($insertStart, $insertStart + $annotation.Length)
$annotationRange.Font.Color = 192   # RGB(0xC0,0x00,0x00) -> C00000

# ---------------------------------------------------------------------------
# 2. Style "Normal (Web)": mark it semi-hidden (Word UI: Style.Hidden ->
#    w:semiHidden). Wrapped defensively since not every host exposes a
#    working setter for this particular property.
# ---------------------------------------------------------------------------

$normalWeb = $d.Styles("Normal (Web)")
try {
    $normalWeb.Hidden = $true
} catch {
    Write-Host "Style.Hidden setter unavailable:" $_
}
